# Edit: "dans la <corr>g</corr>ue<corr>u</corr>le ou aultre endroit. Aprés "
#   ->   "dans la gue<corr>u</corr>le<del><ill/></del> ou aultre endroit. Aprés "
#
# Both before/after read the same once the <tags> are stripped
# ("dans la gueule ou aultre endroit. Aprés "); only the markup
# annotating the handwriting changed: the <corr> that used to wrap the
# first "g" is dropped (its "g" becomes plain "u"... together with the
# surrounding text reflow), the <corr> around the second "u" turns into
# a <del> wrapping an <ill/> (illegible) marker instead of the literal
# "u", and the leading "le" is removed from the run that follows.

$d = $word.ActiveDocument

# Locate the unique anchor "dans la " that starts this run of text.
$anchor = $d.Content
$found = $anchor.Find.Execute("dans la ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "anchor 'dans la ' not found"
}
$base = $anchor.Start

# Donor run: an existing "<m>" run elsewhere in the document already
# carries the exact target formatting for the new "<ill/>" run
# (Courier New / blue 0000ff / sz 18 / szCs 18 / rtl 0).
$donorAnchor = $d.Content
$donorFound = $donorAnchor.Find.Execute("<m>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $donorFound) {
    throw "donor '<m>' run not found"
}
$donorText = $donorAnchor.FormattedText

# Run boundaries, relative to $base, in the *original* document:
#   R1 "dans la "   base+0  .. base+8
#   R2 "<corr>"     base+8  .. base+14   (untouched)
#   R3 "g"          base+14 .. base+15
#   R4 "</corr>"    base+15 .. base+22   (untouched)
#   R5 "ue"         base+22 .. base+24
#   R6 "<corr>"     base+24 .. base+30
#   R7 "u"          base+30 .. base+31
#   R8 "</corr>"    base+31 .. base+38
#   R9 "le ou aultre endroit. Aprés " base+38 .. base+66
#
# Apply edits from the rightmost run back to the leftmost so earlier
# offsets stay valid while later ones shift around.

# R9: drop the leading "le " -> " " (keep " ou aultre endroit. Aprés ")
$r9 = $d.Range($base + 38, $base + 66)
if ($r9.Text -ne "le ou aultre endroit. Aprés ") {
    throw "R9 text mismatch: [$($r9.Text)]"
}
$r9.Text = " ou aultre endroit. Aprés "

# R8: "</corr>" -> "</del>" (formatting unchanged)
$r8 = $d.Range($base + 31, $base + 38)
if ($r8.Text -ne "</corr>") {
    throw "R8 text mismatch: [$($r8.Text)]"
}
$r8.Text = "</del>"

# R7: "u" -> "<ill/>" with new formatting (Courier New, blue, sz 18/18)
$r7 = $d.Range($base + 30, $base + 31)
if ($r7.Text -ne "u") {
    throw "R7 text mismatch: [$($r7.Text)]"
}
$r7.FormattedText = $donorText
$r7b = $d.Range($base + 30, $base + 33)
$r7b.Text = "<ill/>"

# R6: "<corr>" -> "<del>" (formatting unchanged)
$r6 = $d.Range($base + 24, $base + 30)
if ($r6.Text -ne "<corr>") {
    throw "R6 text mismatch: [$($r6.Text)]"
}
$r6.Text = "<del>"

# R5: "ue" -> "le"
$r5 = $d.Range($base + 22, $base + 24)
if ($r5.Text -ne "ue") {
    throw "R5 text mismatch: [$($r5.Text)]"
}
$r5.Text = "le"

# R3: "g" -> "u"
$r3 = $d.Range($base + 14, $base + 15)
if ($r3.Text -ne "g") {
    throw "R3 text mismatch: [$($r3.Text)]"
}
$r3.Text = "u"

# R1: "dans la " -> "dans la gue"
$r1 = $d.Range($base + 0, $base + 8)
if ($r1.Text -ne "dans la ") {
    throw "R1 text mismatch: [$($r1.Text)]"
}
$r1.Text = "dans la gue"

Write-Output "done"
